# This script applies a re-sort / re-numbering of species-observation rows
# (rows 2,3,4,6,7,8,9 on the active sheet). Row 5 is left untouched. The
# edit is effectively a permutation of whole-row content: each target row
# ends up holding the data that used to live in a different row. Columns
# that are identical across every observation row (C,P,S,T,U,V,W,Y,Z,AA,
# AB,AD,AE,AG,AW,AX) do not need to be touched - only columns
# A,B,D,E,F,G,H,I,M,Q,R ever differ between rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData($Row, $A, $B, $D, $E, $F, $G, $H, $Q, $R) {
    $ws.Cells.Item($Row, 1).Value = $A    # A - Id
    $ws.Cells.Item($Row, 2).Value = $B    # B - Taxonsorteringsordning
    $ws.Cells.Item($Row, 4).Value = $D    # D - Rödlistade
    $ws.Cells.Item($Row, 5).Value = $E    # E - TaxonId
    $ws.Cells.Item($Row, 6).Value = $F    # F - Artnamn
    $ws.Cells.Item($Row, 7).Value = $G    # G - Vetenskapligt namn
    $ws.Cells.Item($Row, 8).Value = $H    # H - Auktor
    $ws.Cells.Item($Row, 17).Value = $Q   # Q - Ost
    $ws.Cells.Item($Row, 18).Value = $R   # R - Nord
}

# Row 2 now holds what used to be row 6's data (Skrovellav)
Set-RowData 2 111739317 78579 "NT" 2081 "Skrovellav" "Lobaria scrobiculata" "(Scop.) DC." 573911.5177193542 7172648.020174325

# Row 3 now holds what used to be row 2's data (Lunglav)
Set-RowData 3 111739316 78578 "NT" 6458 "Lunglav" "Lobaria pulmonaria" "(L.) Hoffm." 573904.5013778479 7172636.708955797

# Row 4 now holds what used to be row 8's data (Korallblylav)
Set-RowData 4 111739309 78536 "LC" 229497 "Korallblylav" "Parmeliella triptophylla" "(Ach.) Müll.Arg." 574011.1276117128 7172434.078971106

# Row 6 now holds what used to be row 4's data (Stuplav)
Set-RowData 6 111739315 78605 "LC" 6462 "Stuplav" "Nephroma bellum" "(Spreng.) Tuck." 573904.5013778479 7172636.708955797

# Row 7 now holds what used to be row 9's data (Talltita)
Set-RowData 7 111739307 56543 "NT" 103021 "Talltita" "Poecile montanus" "(Conrad von Baldenstein, 1827)" 573960.5743707293 7172501.399265604

# Row 7 also gains the "Antal" (I) and "Aktivitet" (M) values that used to
# be on row 9. Force text storage for I7 so the numeric-looking "3" is not
# reinterpreted as a number (matches the source cell's text/inlineStr type).
$ws.Cells.Item(7, 9).NumberFormat = "@"
$ws.Cells.Item(7, 9).Value = "3"
$ws.Cells.Item(7, 13).Value = "födosökande"

# Row 8 now holds what used to be row 7's data (Rödbrun blekspik)
Set-RowData 8 111739313 73701 "NT" 1467 "Rödbrun blekspik" "Sclerophora coniophaea" "(Norman) J.Mattsson & Middelb." 574025.0565134182 7172443.417908707

# Row 9 now holds what used to be row 3's data (Garnlav)
Set-RowData 9 111739311 77515 "NT" 6425 "Garnlav" "Alectoria sarmentosa" "(Ach.) Ach." 574011.8892867711 7172473.089384713

# Row 9 no longer carries an Antal/Aktivitet value (that data moved to row 7).
$ws.Cells.Item(9, 9).ClearContents()
$ws.Cells.Item(9, 13).ClearContents()
